$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts existing rows 3..23 down to 4..24)
$ws.Rows.Item(3).Insert()

# Set the new cell's value
$ws.Range("A3").Value = "Navbar - Menü für eingeloogten User"

# Update the selected cell to match the saved view state
$ws.Range("E18").Select()
